# Bump the auto-updating "date" footer placeholder text from 18/07/2018 to
# 19/07/2018 everywhere it is cached: on the Slide Master and on every
# Slide Layout (the diff does not touch the individual slides themselves).

$p = $ppt.ActivePresentation

$oldDate = "18/07/2018"
$newDate = "19/07/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout hanging off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
